$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52; existing rows 52:57 shift down to 53:58
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with the new record's data
$ws.Cells.Item(52, 1).Value = 5
$ws.Cells.Item(52, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(52, 3).Value = "Maule"
$ws.Cells.Item(52, 4).Value = 44474
$ws.Cells.Item(52, 5).Value = 7
$ws.Cells.Item(52, 6).Value = 100112013
$ws.Cells.Item(52, 7).Value = "Alcachofa"
$ws.Cells.Item(52, 8).Value = "Madrigal"
$ws.Cells.Item(52, 9).Value = "Primera"
$ws.Cells.Item(52, 10).Value = 500
$ws.Cells.Item(52, 11).Value = 10000
$ws.Cells.Item(52, 12).Value = 10000
$ws.Cells.Item(52, 13).Value = 10000
$ws.Cells.Item(52, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(52, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(52, 16).Value = 250
$ws.Cells.Item(52, 17).Value = 40
$ws.Cells.Item(52, 18).Value = "Hortaliza"
